$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 2954
$ws.Cells.Item(32, 10).Value = 3066.8
$ws.Cells.Item(32, 12).Value = 3066.8
$ws.Cells.Item(32, 14).Value = -3718.8
$ws.Cells.Item(53, 8).Value = 314
$ws.Cells.Item(53, 9).Value = 213.6
$ws.Cells.Item(53, 11).Value = 213.6
$ws.Cells.Item(53, 13).Value = 423.4
$ws.Cells.Item(107, 8).Value = 933.5
$ws.Cells.Item(107, 9).Value = 391.66666
$ws.Cells.Item(107, 11).Value = 391.66666
$ws.Cells.Item(107, 13).Value = 1528.33334
$ws.Cells.Item(112, 8).Value = 3125.2354
$ws.Cells.Item(112, 9).Value = 1394
$ws.Cells.Item(112, 10).Value = 3496.2144
$ws.Cells.Item(112, 11).Value = 4182
$ws.Cells.Item(112, 12).Value = 10488.6432
$ws.Cells.Item(112, 13).Value = -3074
$ws.Cells.Item(112, 14).Value = -12704.6432
$ws.Cells.Item(132, 8).Value = 2979.4443
$ws.Cells.Item(132, 9).Value = 2726.875
$ws.Cells.Item(132, 10).Value = 5000
$ws.Cells.Item(132, 11).Value = 8180.625
$ws.Cells.Item(132, 12).Value = 15000
$ws.Cells.Item(132, 13).Value = -5650.625
$ws.Cells.Item(132, 14).Value = -20060

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(39, 8).Value = 1758.5
$ws.Cells.Item(39, 9).Value = 972
$ws.Cells.Item(39, 10).Value = 4118
$ws.Cells.Item(39, 11).Value = 972
$ws.Cells.Item(39, 12).Value = 4118
$ws.Cells.Item(39, 13).Value = -452
$ws.Cells.Item(39, 14).Value = -5158
$ws.Cells.Item(46, 8).Value = 9573
$ws.Cells.Item(46, 10).Value = 9514
$ws.Cells.Item(46, 12).Value = 9514
$ws.Cells.Item(46, 14).Value = -10152
$ws.Cells.Item(121, 8).Value = 52500
$ws.Cells.Item(121, 9).Value = 0
$ws.Cells.Item(121, 10).Value = 52500
$ws.Cells.Item(121, 11).Value = 0
$ws.Cells.Item(121, 12).Value = 52500
$ws.Cells.Item(121, 14).Value = -55994
$ws.Cells.Item(122, 8).Value = 20380.438
$ws.Cells.Item(122, 9).Value = 20380.438
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 61141.314
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -58691.314
$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 11).Value = 0
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(126, 8).Value = 0
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(127, 8).Value = 0
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 11).Value = 0
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(128, 8).Value = 0
$ws.Cells.Item(128, 9).Value = 0
$ws.Cells.Item(128, 10).Value = 0
$ws.Cells.Item(128, 11).Value = 0
$ws.Cells.Item(128, 12).Value = 0
$ws.Cells.Item(129, 8).Value = 0
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 0
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 0
$ws.Cells.Item(130, 8).Value = 0
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 0
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 12).Value = 0
$ws.Cells.Item(131, 8).Value = 0
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(132, 8).Value = 1956.3846
$ws.Cells.Item(132, 9).Value = 1917.25
$ws.Cells.Item(132, 10).Value = 2019
$ws.Cells.Item(132, 11).Value = 5751.75
$ws.Cells.Item(132, 12).Value = 6057
$ws.Cells.Item(132, 13).Value = -3221.75
$ws.Cells.Item(132, 14).Value = -11117
$ws.Cells.Item(133, 8).Value = 147443.75
$ws.Cells.Item(133, 9).Value = 144000
$ws.Cells.Item(133, 10).Value = 148591.67
$ws.Cells.Item(133, 11).Value = 144000
$ws.Cells.Item(133, 12).Value = 148591.67
$ws.Cells.Item(133, 13).Value = -141470
$ws.Cells.Item(133, 14).Value = -153651.67
$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(135, 8).Value = 387976.34
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 387976.34
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 387976.34
$ws.Cells.Item(135, 14).Value = -398116.34
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(138, 8).Value = 250000
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 250000
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 250000
$ws.Cells.Item(138, 14).Value = -260280
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 0

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(37, 8).Value = 14816.333
$ws.Cells.Item(37, 9).Value = 6224.5
$ws.Cells.Item(37, 11).Value = 6224.5
$ws.Cells.Item(37, 13).Value = -6087.5

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1950.5
$ws.Cells.Item(22, 9).Value = 1998
$ws.Cells.Item(22, 11).Value = 1998
$ws.Cells.Item(22, 13).Value = -1648
$ws.Cells.Item(31, 8).Value = 1104.88
$ws.Cells.Item(31, 10).Value = 1314
$ws.Cells.Item(31, 12).Value = 1314
$ws.Cells.Item(31, 14).Value = -1904
$ws.Cells.Item(34, 8).Value = 1104.88
$ws.Cells.Item(34, 10).Value = 1314
$ws.Cells.Item(34, 12).Value = 1314
$ws.Cells.Item(34, 14).Value = -1718
$ws.Cells.Item(51, 8).Value = 50000
$ws.Cells.Item(51, 10).Value = 70000
$ws.Cells.Item(51, 12).Value = 70000
$ws.Cells.Item(51, 14).Value = -71472
$ws.Cells.Item(58, 8).Value = 2796.125
$ws.Cells.Item(58, 9).Value = 2084.5715
$ws.Cells.Item(58, 11).Value = 2084.5715
$ws.Cells.Item(58, 13).Value = -1881.5715
$ws.Cells.Item(61, 8).Value = 50000
$ws.Cells.Item(61, 10).Value = 70000
$ws.Cells.Item(61, 12).Value = 70000
$ws.Cells.Item(61, 14).Value = -70696
$ws.Cells.Item(129, 8).Value = 0
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 0
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 0
$ws.Cells.Item(130, 8).Value = 89999
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 89999
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 12).Value = 89999
$ws.Cells.Item(130, 14).Value = -100039
$ws.Cells.Item(131, 8).Value = 49878.332
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 49878.332
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 12).Value = 49878.332
$ws.Cells.Item(131, 14).Value = -59958.332
$ws.Cells.Item(132, 8).Value = 7877.3335
$ws.Cells.Item(132, 9).Value = 7893.5713
$ws.Cells.Item(132, 10).Value = 7650
$ws.Cells.Item(132, 11).Value = 23680.7139
$ws.Cells.Item(132, 12).Value = 22950
$ws.Cells.Item(132, 13).Value = -21150.7139
$ws.Cells.Item(132, 14).Value = -28010
$ws.Cells.Item(133, 8).Value = 63999
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 63999
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 63999
$ws.Cells.Item(133, 14).Value = -69059
$ws.Cells.Item(134, 8).Value = 2702.6667
$ws.Cells.Item(134, 9).Value = 2702.6667
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 8108.000100000001
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).Value = -5573.000100000001
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(138, 8).Value = 60000
$ws.Cells.Item(138, 9).Value = 60000
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 11).Value = 60000
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 13).Value = -54860
$ws.Cells.Item(139, 8).Value = 40000
$ws.Cells.Item(139, 9).Value = 40000
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 11).Value = 40000
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 13).Value = -34860
$ws.Cells.Item(140, 8).Value = 80000
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 80000
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 80000
$ws.Cells.Item(140, 14).Value = -90360
$ws.Cells.Item(141, 8).Value = 240200
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 240200
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 240200
$ws.Cells.Item(141, 14).Value = -250560
$ws.Cells.Item(136, 8).Value = 2796.125
$ws.Cells.Item(136, 9).Value = 2084.5715
$ws.Cells.Item(136, 11).Value = 6253.7145
$ws.Cells.Item(136, 13).Value = -3703.7145

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 48.588234
$ws.Cells.Item(2, 9).Value = 26.6
$ws.Cells.Item(2, 10).Value = 80
$ws.Cells.Item(2, 11).Value = 26.6
$ws.Cells.Item(2, 12).Value = 80
$ws.Cells.Item(2, 13).Value = 86.4
$ws.Cells.Item(2, 14).Value = -306
$ws.Cells.Item(3, 8).Value = 650633.9
$ws.Cells.Item(3, 9).Value = 917225
$ws.Cells.Item(3, 10).Value = 250747.25
$ws.Cells.Item(3, 11).Value = 917225
$ws.Cells.Item(3, 12).Value = 250747.25
$ws.Cells.Item(3, 13).Value = -917109
$ws.Cells.Item(3, 14).Value = -250979.25
$ws.Cells.Item(9, 8).Value = 530
$ws.Cells.Item(9, 9).Value = 595
$ws.Cells.Item(9, 10).Value = 400
$ws.Cells.Item(9, 11).Value = 595
$ws.Cells.Item(9, 12).Value = 400
$ws.Cells.Item(9, 13).Value = -425
$ws.Cells.Item(9, 14).Value = -740
$ws.Cells.Item(10, 8).Value = 7476.25
$ws.Cells.Item(10, 9).Value = 8300.333
$ws.Cells.Item(10, 10).Value = 5004
$ws.Cells.Item(10, 11).Value = 8300.333
$ws.Cells.Item(10, 12).Value = 5004
$ws.Cells.Item(10, 13).Value = -8131.333000000001
$ws.Cells.Item(10, 14).Value = -5342
$ws.Cells.Item(11, 8).Value = 2063750.9
$ws.Cells.Item(11, 9).Value = 3100000
$ws.Cells.Item(11, 10).Value = 336669
$ws.Cells.Item(11, 11).Value = 3100000
$ws.Cells.Item(11, 12).Value = 336669
$ws.Cells.Item(11, 13).Value = -3099861
$ws.Cells.Item(11, 14).Value = -336947
$ws.Cells.Item(12, 8).Value = 3000
$ws.Cells.Item(12, 9).Value = 3000
$ws.Cells.Item(12, 11).Value = 3000
$ws.Cells.Item(12, 13).Value = -2860
$ws.Cells.Item(13, 8).Value = 240.66667
$ws.Cells.Item(13, 10).Value = 348.5
$ws.Cells.Item(13, 12).Value = 348.5
$ws.Cells.Item(13, 14).Value = -626.5
$ws.Cells.Item(63, 8).Value = 49999
$ws.Cells.Item(63, 9).Value = 49999
$ws.Cells.Item(63, 10).Value = 49999
$ws.Cells.Item(63, 11).Value = 49999
$ws.Cells.Item(63, 12).Value = 49999
$ws.Cells.Item(63, 13).Value = -49313
$ws.Cells.Item(63, 14).Value = -51371
$ws.Cells.Item(66, 8).Value = 49999
$ws.Cells.Item(66, 9).Value = 49999
$ws.Cells.Item(66, 10).Value = 49999
$ws.Cells.Item(66, 11).Value = 149997
$ws.Cells.Item(66, 12).Value = 149997
$ws.Cells.Item(66, 13).Value = -146565
$ws.Cells.Item(66, 14).Value = -156861
$ws.Cells.Item(70, 8).Value = 4002
$ws.Cells.Item(70, 9).Value = 3998.5
$ws.Cells.Item(70, 10).Value = 4009
$ws.Cells.Item(70, 11).Value = 3998.5
$ws.Cells.Item(70, 12).Value = 4009
$ws.Cells.Item(70, 13).Value = -3728.5
$ws.Cells.Item(70, 14).Value = -4549
$ws.Cells.Item(73, 8).Value = 4002
$ws.Cells.Item(73, 9).Value = 3998.5
$ws.Cells.Item(73, 10).Value = 4009
$ws.Cells.Item(73, 11).Value = 3998.5
$ws.Cells.Item(73, 12).Value = 4009
$ws.Cells.Item(73, 13).Value = -3062.5
$ws.Cells.Item(73, 14).Value = -5881
$ws.Cells.Item(125, 8).Value = 80000
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 10).Value = 80000
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 12).Value = 80000
$ws.Cells.Item(125, 14).Value = -84920
$ws.Cells.Item(126, 8).Value = 5449
$ws.Cells.Item(126, 9).Value = 4231
$ws.Cells.Item(126, 10).Value = 5719.6665
$ws.Cells.Item(126, 11).Value = 12693
$ws.Cells.Item(126, 12).Value = 17158.9995
$ws.Cells.Item(126, 13).Value = -10223
$ws.Cells.Item(126, 14).Value = -22098.9995
$ws.Cells.Item(127, 8).Value = 0
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 11).Value = 0
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(128, 8).Value = 0
$ws.Cells.Item(128, 9).Value = 0
$ws.Cells.Item(128, 10).Value = 0
$ws.Cells.Item(128, 11).Value = 0
$ws.Cells.Item(128, 12).Value = 0
$ws.Cells.Item(129, 8).Value = 116519.336
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 116519.336
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 116519.336
$ws.Cells.Item(129, 14).Value = -126519.336
$ws.Cells.Item(130, 8).Value = 89999
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 89999
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 12).Value = 89999
$ws.Cells.Item(130, 14).Value = -100039
$ws.Cells.Item(131, 8).Value = 0
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(132, 8).Value = 9502
$ws.Cells.Item(132, 9).Value = 9502
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 28506
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -25976
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(134, 8).Value = 46666.332
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 10).Value = 46666.332
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 12).Value = 139998.996
$ws.Cells.Item(134, 14).Value = -145068.996
$ws.Cells.Item(135, 8).Value = 255000
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 255000
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 255000
$ws.Cells.Item(135, 14).Value = -265140
$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 0

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(11, 8).Value = 14000
$ws.Cells.Item(11, 10).Value = 14000
$ws.Cells.Item(11, 12).Value = 14000
$ws.Cells.Item(11, 14).Value = -14280
$ws.Cells.Item(39, 8).Value = 19999
$ws.Cells.Item(39, 10).Value = 19999
$ws.Cells.Item(39, 12).Value = 19999
$ws.Cells.Item(39, 14).Value = -20919
$ws.Cells.Item(55, 8).Value = 1197.4667
$ws.Cells.Item(55, 9).Value = 897.7
$ws.Cells.Item(55, 11).Value = 897.7
$ws.Cells.Item(55, 13).Value = -724.7
$ws.Cells.Item(122, 8).Value = 5868.8335
$ws.Cells.Item(122, 9).Value = 4719.952
$ws.Cells.Item(122, 11).Value = 14159.856
$ws.Cells.Item(122, 13).Value = -11709.856

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1627.6666
$ws.Cells.Item(107, 9).Value = 1512.7273
$ws.Cells.Item(107, 11).Value = 4538.1819
$ws.Cells.Item(107, 13).Value = -2618.1819
$ws.Cells.Item(113, 8).Value = 431.66666
$ws.Cells.Item(113, 9).Value = 397.66666
$ws.Cells.Item(113, 11).Value = 1192.99998
$ws.Cells.Item(113, 13).Value = 977.00002
$ws.Cells.Item(122, 8).Value = 3859.9443
$ws.Cells.Item(122, 9).Value = 2546.6924
$ws.Cells.Item(122, 10).Value = 7274.4
$ws.Cells.Item(122, 11).Value = 7640.0772
$ws.Cells.Item(122, 12).Value = 21823.2
$ws.Cells.Item(122, 13).Value = -5190.0772
$ws.Cells.Item(122, 14).Value = -26723.2
$ws.Cells.Item(132, 8).Value = 1511.12
$ws.Cells.Item(132, 9).Value = 1526.6875
$ws.Cells.Item(132, 10).Value = 1483.4445
$ws.Cells.Item(132, 11).Value = 4580.0625
$ws.Cells.Item(132, 12).Value = 4450.333500000001
$ws.Cells.Item(132, 13).Value = -2050.0625
$ws.Cells.Item(132, 14).Value = -9510.3335
$ws.Cells.Item(136, 8).Value = 3349.087
$ws.Cells.Item(136, 9).Value = 3282.1
$ws.Cells.Item(136, 11).Value = 9846.3
$ws.Cells.Item(136, 13).Value = -7296.299999999999
